$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.901461666666667
$ws.Range("H2").Value = 14.704385
$ws.Range("I2").Value = 0.2124427850531459
$ws.Range("J2").Value = 0.2124427850531459
$ws.Range("M2").Value = 2.407428333333334
$ws.Range("N2").Value = 7.222285000000001
$ws.Range("O2").Value = 0.3774919259534536
$ws.Range("P2").Value = 0.3774919259534536
$ws.Range("Q2").Value = 11.79991769108056
$ws.Range("R2").Value = 106.199259219725
$ws.Range("S2").Value = 0.08019543608462761
$ws.Range("T2").Value = 0.08019543608462761
$ws.Range("G3").Value = 4.901461666666667
$ws.Range("H3").Value = 14.704385
$ws.Range("I3").Value = 0.2124427850531459
$ws.Range("J3").Value = 0.2124427850531459
$ws.Range("O3").Value = 0.3035072644205163
$ws.Range("P3").Value = 0.3035072644205163
$ws.Range("Q3").Value = 9.487251229973891
$ws.Range("R3").Value = 85.38526106976501
$ws.Range("S3").Value = 0.06447792853735605
$ws.Range("T3").Value = 0.06447792853735605
$ws.Range("G4").Value = 4.901461666666667
$ws.Range("H4").Value = 14.704385
$ws.Range("I4").Value = 0.2124427850531459
$ws.Range("J4").Value = 0.2124427850531459
$ws.Range("M4").Value = 1.020054666666667
$ws.Range("N4").Value = 3.060164
$ws.Range("O4").Value = 0.1599476068991219
$ws.Range("P4").Value = 0.1599476068991219
$ws.Range("Q4").Value = 4.999758846571112
$ws.Range("R4").Value = 44.99782961914001
$ws.Range("S4").Value = 0.03397971507223522
$ws.Range("T4").Value = 0.03397971507223522
$ws.Range("G5").Value = 4.901461666666667
$ws.Range("H5").Value = 14.704385
$ws.Range("I5").Value = 0.2124427850531459
$ws.Range("J5").Value = 0.2124427850531459
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.014350666666667
$ws.Range("N5").Value = 3.043052
$ws.Range("O5").Value = 0.1590532027269083
$ws.Range("P5").Value = 0.1590532027269083
$ws.Range("Q5").Value = 4.971800909224444
$ws.Range("R5").Value = 44.74620818302
$ws.Range("S5").Value = 0.03378970535892701
$ws.Range("T5").Value = 0.03378970535892701
$ws.Range("I6").Value = 0.1372144215401173
$ws.Range("J6").Value = 0.1372144215401173
$ws.Range("M6").Value = 2.407428333333334
$ws.Range("N6").Value = 7.222285000000001
$ws.Range("O6").Value = 0.3774919259534536
$ws.Range("P6").Value = 0.3774919259534536
$ws.Range("Q6").Value = 7.621435012714446
$ws.Range("R6").Value = 68.59291511443001
$ws.Range("S6").Value = 0.05179733625576793
$ws.Range("T6").Value = 0.05179733625576793
$ws.Range("I7").Value = 0.1372144215401173
$ws.Range("J7").Value = 0.1372144215401173
$ws.Range("O7").Value = 0.3035072644205163
$ws.Range("P7").Value = 0.3035072644205163
$ws.Range("S7").Value = 0.04164557372068457
$ws.Range("T7").Value = 0.04164557372068457
$ws.Range("I8").Value = 0.1372144215401173
$ws.Range("J8").Value = 0.1372144215401173
$ws.Range("M8").Value = 1.020054666666667
$ws.Range("N8").Value = 3.060164
$ws.Range("O8").Value = 0.1599476068991219
$ws.Range("P8").Value = 0.1599476068991219
$ws.Range("Q8").Value = 3.22928838369689
$ws.Range("R8").Value = 29.063595453272
$ws.Range("S8").Value = 0.02194711835738908
$ws.Range("T8").Value = 0.02194711835738908
$ws.Range("I9").Value = 0.1372144215401173
$ws.Range("J9").Value = 0.1372144215401173
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 1.014350666666667
$ws.Range("N9").Value = 3.043052
$ws.Range("O9").Value = 0.1590532027269083
$ws.Range("P9").Value = 0.1590532027269083
$ws.Range("Q9").Value = 3.211230664299556
$ws.Range("R9").Value = 28.901075978696
$ws.Range("S9").Value = 0.02182439320627572
$ws.Range("T9").Value = 0.02182439320627572
$ws.Range("G10").Value = 1.206743666666667
$ws.Range("H10").Value = 3.620231
$ws.Range("I10").Value = 0.05230357857032003
$ws.Range("J10").Value = 0.05230357857032004
$ws.Range("M10").Value = 2.407428333333334
$ws.Range("N10").Value = 7.222285000000001
$ws.Range("O10").Value = 0.3774919259534536
$ws.Range("P10").Value = 0.3774919259534536
$ws.Range("Q10").Value = 2.905148894203889
$ws.Range("R10").Value = 26.146340047835
$ws.Range("S10").Value = 0.01974417860876789
$ws.Range("T10").Value = 0.01974417860876789
$ws.Range("G11").Value = 1.206743666666667
$ws.Range("H11").Value = 3.620231
$ws.Range("I11").Value = 0.05230357857032003
$ws.Range("J11").Value = 0.05230357857032004
$ws.Range("O11").Value = 0.3035072644205163
$ws.Range("P11").Value = 0.3035072644205163
$ws.Range("Q11").Value = 2.335768616473223
$ws.Range("R11").Value = 21.021917548259
$ws.Range("S11").Value = 0.01587451605128137
$ws.Range("T11").Value = 0.01587451605128137
$ws.Range("G12").Value = 1.206743666666667
$ws.Range("H12").Value = 3.620231
$ws.Range("I12").Value = 0.05230357857032003
$ws.Range("J12").Value = 0.05230357857032004
$ws.Range("M12").Value = 1.020054666666667
$ws.Range("N12").Value = 3.060164
$ws.Range("O12").Value = 0.1599476068991219
$ws.Range("P12").Value = 0.1599476068991219
$ws.Range("Q12").Value = 1.230944508653778
$ws.Range("R12").Value = 11.078500577884
$ws.Range("S12").Value = 0.008365832224582882
$ws.Range("T12").Value = 0.008365832224582884
$ws.Range("G13").Value = 1.206743666666667
$ws.Range("H13").Value = 3.620231
$ws.Range("I13").Value = 0.05230357857032003
$ws.Range("J13").Value = 0.05230357857032004
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 1.014350666666667
$ws.Range("N13").Value = 3.043052
$ws.Range("O13").Value = 0.1590532027269083
$ws.Range("P13").Value = 0.1590532027269083
$ws.Range("Q13").Value = 1.224061242779111
$ws.Range("R13").Value = 11.016551185012
$ws.Range("S13").Value = 0.008319051685687886
$ws.Range("T13").Value = 0.008319051685687888
$ws.Range("G14").Value = 13.79790933333333
$ws.Range("H14").Value = 41.393728
$ws.Range("I14").Value = 0.5980392148364168
$ws.Range("J14").Value = 0.5980392148364169
$ws.Range("M14").Value = 2.407428333333334
$ws.Range("N14").Value = 7.222285000000001
$ws.Range("O14").Value = 0.3774919259534536
$ws.Range("P14").Value = 0.3774919259534536
$ws.Range("Q14").Value = 33.21747786983111
$ws.Range("R14").Value = 298.95730082848
$ws.Range("S14").Value = 0.2257549750042902
$ws.Range("T14").Value = 0.2257549750042902
$ws.Range("G15").Value = 13.79790933333333
$ws.Range("H15").Value = 41.393728
$ws.Range("I15").Value = 0.5980392148364168
$ws.Range("J15").Value = 0.5980392148364169
$ws.Range("O15").Value = 0.3035072644205163
$ws.Range("P15").Value = 0.3035072644205163
$ws.Range("Q15").Value = 26.70718271326578
$ws.Range("R15").Value = 240.364644419392
$ws.Range("S15").Value = 0.1815092461111943
$ws.Range("T15").Value = 0.1815092461111943
$ws.Range("G16").Value = 13.79790933333333
$ws.Range("H16").Value = 41.393728
$ws.Range("I16").Value = 0.5980392148364168
$ws.Range("J16").Value = 0.5980392148364169
$ws.Range("M16").Value = 1.020054666666667
$ws.Range("N16").Value = 3.060164
$ws.Range("O16").Value = 0.1599476068991219
$ws.Range("P16").Value = 0.1599476068991219
$ws.Range("Q16").Value = 14.07462180571022
$ws.Range("R16").Value = 126.671596251392
$ws.Range("S16").Value = 0.09565494124491467
$ws.Range("T16").Value = 0.09565494124491468
$ws.Range("G17").Value = 13.79790933333333
$ws.Range("H17").Value = 41.393728
$ws.Range("I17").Value = 0.5980392148364168
$ws.Range("J17").Value = 0.5980392148364169
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 1.014350666666667
$ws.Range("N17").Value = 3.043052
$ws.Range("O17").Value = 0.1590532027269083
$ws.Range("P17").Value = 0.1590532027269083
$ws.Range("Q17").Value = 13.99591853087289
$ws.Range("R17").Value = 125.963266777856
$ws.Range("S17").Value = 0.09512005247601764
$ws.Range("T17").Value = 0.09512005247601765
